$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-23 Tuesday" "2025-12-24 Wednesday"

Replace-Text "19×44=" "73×41="
Replace-Text "20×70=" "39×71="
Replace-Text "22×75=" "11×28="
Replace-Text "91×18=" "12×83="
Replace-Text "25×42=" "16×13="
Replace-Text "60×21=" "50×25="
Replace-Text "59×81=" "62×51="
Replace-Text "69×42=" "90×65="
Replace-Text "34×61=" "55×41="
Replace-Text "11×72=" "30×38="
Replace-Text "16×53=" "13×48="
Replace-Text "44×74=" "65×61="
Replace-Text "82×46=" "15×78="
Replace-Text "25×21=" "56×85="
Replace-Text "83×26=" "31×61="
Replace-Text "63×98=" "20×77="
Replace-Text "87×96=" "36×82="
Replace-Text "82×65=" "67×47="
Replace-Text "86×23=" "57×85="
Replace-Text "88×65=" "89×17="
Replace-Text "84×41=" "99×59="
Replace-Text "81×30=" "13×59="
Replace-Text "42×51=" "64×17="
Replace-Text "83×99=" "63×85="
Replace-Text "73×92=" "92×89="
